$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Insert the six new paragraphs (a blank N1 spacer, the sealing statement,
# a blank line, the signatory block and a trailing blank N1 spacer) after the
# existing final paragraph of the document, immediately before the sectPr.
#
# NOTE: this runtime has an off-by-one quirk where inserting text at a
# position that is the absolute end of the document content lands the text
# one character earlier than the reported Range.Start/End. To avoid it we
# always create the *next* empty paragraph first (so the paragraph we are
# about to fill text into is no longer the very last position in the
# document) before writing any text into the current paragraph.
# ---------------------------------------------------------------------------

function New-TrailingParagraph {
    $r = $d.Range($d.Content.End, $d.Content.End)
    $r.InsertParagraphAfter()
    return $d.Paragraphs.Last
}

# Paragraph 1: blank spacer, style N1, numbering removed.
$p1 = New-TrailingParagraph

# Paragraph 2: sealing statement (style SigBlock).
$p2 = New-TrailingParagraph

# p1 needs no text - just style + remove numbering.
$p1.Range.Style = "N1"
$p1.Range.ListFormat.RemoveNumbers()

# Paragraph 3: blank line (style SigBlock) - create now so p2's text fill
# below is not at the absolute end of the document.
$p3 = New-TrailingParagraph

# Fill paragraph 2's text.
$p2.Range.Style = "SigBlock"
$p2start = $p2.Range.Start
$p2text = $d.Range($p2start, $p2start)
$p2text.Text = "Sealed with the Official Seal of the Department of Agriculture, Environment and Rural Affairs on XXXX 2025."
$p2text.Style = "Sigsignatory"

# Paragraph 4: "Senior Officer" line (style SigBlock).
$p4 = New-TrailingParagraph

# p3 stays blank - just the style.
$p3.Range.Style = "SigBlock"

# Paragraph 5: title line (style SigBlock).
$p5 = New-TrailingParagraph

# Fill paragraph 4's content: a tab followed by the styled "Senior Officer".
$p4.Range.Style = "SigBlock"
$p4tabStart = $p4.Range.Start
$p4tab = $d.Range($p4tabStart, $p4tabStart)
$p4tab.Text = "`t"
$p4textStart = $p4tab.End
$p4text = $d.Range($p4textStart, $p4textStart)
$p4text.Text = "Senior Officer"
$p4text.Style = "SigSignee"

# Paragraph 6: trailing blank spacer, style N1, numbering removed.
$p6 = New-TrailingParagraph

# Fill paragraph 5's content: a tab followed by the styled title text.
$p5.Range.Style = "SigBlock"
$p5tabStart = $p5.Range.Start
$p5tab = $d.Range($p5tabStart, $p5tabStart)
$p5tab.Text = "`t"
$p5textStart = $p5tab.End
$p5text = $d.Range($p5textStart, $p5textStart)
$p5text.Text = "A senior officer of the Department of Agriculture, Environment and Rural Affairs"
$p5text.Style = "Sigtitle"

# p6 stays blank - style N1 with numbering removed.
$p6.Range.Style = "N1"
$p6.Range.ListFormat.RemoveNumbers()

Write-Host "Done inserting signature block paragraphs."
